$ws = $excel.ActiveWorkbook.ActiveSheet

$changes = @(
    @{Row=7; I='sd'; J='Statement-non-opinion'},
    @{Row=13; I='aa'; J='Agree/Accept'},
    @{Row=22; I='sv'; J='Statement-opinion'},
    @{Row=35; I='sd'; J='Statement-non-opinion'},
    @{Row=38; I='sd'; J='Statement-non-opinion'},
    @{Row=54; I='sd'; J='Statement-non-opinion'},
    @{Row=83; I='aa'; J='Agree/Accept'},
    @{Row=92; I='sv'; J='Statement-opinion'},
    @{Row=93; I='sd'; J='Statement-non-opinion'},
    @{Row=95; I='%'; J='Uninterpretable'},
    @{Row=101; I='sd'; J='Statement-non-opinion'},
    @{Row=114; I='sv'; J='Statement-opinion'},
    @{Row=123; I='aa'; J='Agree/Accept'},
    @{Row=124; I='%'; J='Uninterpretable'},
    @{Row=148; I='sv'; J='Statement-opinion'},
    @{Row=152; I='sd'; J='Statement-non-opinion'},
    @{Row=156; I='ba'; J='Appreciation'},
    @{Row=166; I='sd'; J='Statement-non-opinion'},
    @{Row=173; I='aa'; J='Agree/Accept'},
    @{Row=179; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=200; I='sd'; J='Statement-non-opinion'},
    @{Row=201; I='aa'; J='Agree/Accept'},
    @{Row=206; I='sd'; J='Statement-non-opinion'},
    @{Row=247; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=248; I='sd'; J='Statement-non-opinion'},
    @{Row=253; I='aa'; J='Agree/Accept'},
    @{Row=256; I='ba'; J='Appreciation'},
    @{Row=258; I='ba'; J='Appreciation'},
    @{Row=272; I='aa'; J='Agree/Accept'},
    @{Row=282; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=286; I='aa'; J='Agree/Accept'},
    @{Row=289; I='qy'; J='Yes-No-Question'},
    @{Row=301; I='sv'; J='Statement-opinion'},
    @{Row=303; I='aa'; J='Agree/Accept'},
    @{Row=314; I='aa'; J='Agree/Accept'},
    @{Row=324; I='sd'; J='Statement-non-opinion'},
    @{Row=333; I='sd'; J='Statement-non-opinion'},
    @{Row=334; I='sd'; J='Statement-non-opinion'},
    @{Row=337; I='sd'; J='Statement-non-opinion'},
    @{Row=357; I='aa'; J='Agree/Accept'},
    @{Row=364; I='sd'; J='Statement-non-opinion'},
    @{Row=373; I='sd'; J='Statement-non-opinion'},
    @{Row=376; I='%'; J='Uninterpretable'},
    @{Row=381; I='sv'; J='Statement-opinion'},
    @{Row=396; I='sd'; J='Statement-non-opinion'},
    @{Row=407; I='sv'; J='Statement-opinion'},
    @{Row=412; I='sv'; J='Statement-opinion'},
    @{Row=413; I='aa'; J='Agree/Accept'},
    @{Row=414; I='ba'; J='Appreciation'},
    @{Row=419; I='sv'; J='Statement-opinion'},
    @{Row=425; I='sd'; J='Statement-non-opinion'},
    @{Row=428; I='sd'; J='Statement-non-opinion'},
    @{Row=435; I='sv'; J='Statement-opinion'},
    @{Row=436; I='sd'; J='Statement-non-opinion'},
    @{Row=449; I='sd'; J='Statement-non-opinion'},
    @{Row=470; I='aa'; J='Agree/Accept'},
    @{Row=478; I='sd'; J='Statement-non-opinion'},
    @{Row=490; I='aa'; J='Agree/Accept'},
    @{Row=491; I='ba'; J='Appreciation'},
    @{Row=510; I='%'; J='Uninterpretable'},
    @{Row=518; I='aa'; J='Agree/Accept'},
    @{Row=524; I='aa'; J='Agree/Accept'},
    @{Row=525; I='sd'; J='Statement-non-opinion'},
    @{Row=526; I='sd'; J='Statement-non-opinion'},
    @{Row=531; I='aa'; J='Agree/Accept'},
    @{Row=546; I='aa'; J='Agree/Accept'},
    @{Row=551; I='sd'; J='Statement-non-opinion'},
    @{Row=553; I='sd'; J='Statement-non-opinion'},
    @{Row=554; I='sd'; J='Statement-non-opinion'},
    @{Row=558; I='sd'; J='Statement-non-opinion'},
    @{Row=573; I='ba'; J='Appreciation'},
    @{Row=581; I='sd'; J='Statement-non-opinion'},
    @{Row=582; I='sd'; J='Statement-non-opinion'},
    @{Row=585; I='aa'; J='Agree/Accept'},
    @{Row=588; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=589; I='ba'; J='Appreciation'},
    @{Row=618; I='sd'; J='Statement-non-opinion'},
    @{Row=635; I='aa'; J='Agree/Accept'},
    @{Row=645; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=646; I='sd'; J='Statement-non-opinion'},
    @{Row=651; I='sd'; J='Statement-non-opinion'},
    @{Row=661; I='sd'; J='Statement-non-opinion'},
    @{Row=695; I='sd'; J='Statement-non-opinion'},
    @{Row=714; I='ba'; J='Appreciation'},
    @{Row=725; I='b'; J='Acknowledge (Backchannel)'},
    @{Row=734; I='ba'; J='Appreciation'},
    @{Row=744; I='fc'; J='Conventional-closing'}
)

foreach ($chg in $changes) {
    $ws.Cells.Item($chg.Row, 9).Value = $chg.I
    $ws.Cells.Item($chg.Row, 10).Value = $chg.J
}
